$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '31.191.87'
$ws.Range("E2").Value = '  +2.73%  '

$ws.Range("D3").Value = '1.984.96'
$ws.Range("E3").Value = '  +6.14%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9994'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.8072'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +71.47%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '253.37'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.95%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9982'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.20%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3448'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +19.86%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.61'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +16.38%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06954'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.79%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8410'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +16.07%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08109'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.27%  '

$ws.Range("D13").Value = '1.984.77'
$ws.Range("E13").Value = '  +6.08%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '100.42'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.45%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.502'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +7.33%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '271.14'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.82%  '

$ws.Range("D17").Value = '31.211.73'
$ws.Range("E17").Value = '  +2.82%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.93'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +7.23%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007953'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.05%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.797'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +10.73%  '

$ws.Range("D21").Value = '2.245.35'
$ws.Range("E21").Value = '  +6.60%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9979'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.21%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9991'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.09%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.924'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +11.15%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.722'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +7.50%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1491'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +54.64%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '164.25'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.50%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.92'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.68%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.186'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +16.54%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.566'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.36%  '

$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.353'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.43%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.564'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +8.37%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.308'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.84%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05158'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.36%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.215'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +8.49%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7571'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +10.03%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.766'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.19%  '

$ws.Range("E38").Value = '  +6.27%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.909'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.55%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.587'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.86%  '

$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '78.04'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.28%  '

$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4681'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +10.91%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.065'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.74%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8525'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.45%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '104.47'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.57%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9977'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.16%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.953'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.81%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.496'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +7.79%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.72'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.20%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4293'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +9.73%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1192'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +12.44%  '
